# Generate Report for Archive
# - Status changes from "Ready for handoff" to "In Translation" for the
#   single tracked file, reflected on the Overview sheet (zh-cn + de-de
#   status columns) and on each per-locale sheet's own Status column.
# - The Status column narrows afterwards (the new text is shorter), so the
#   affected columns are resized to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Overview sheet: column E = zh-cn status, column F = de-de status
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C = Status
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C = Status
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
